$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new bowler records to the table (rows 8 and 9)
$row8 = @("XC Bartlett", "2020-2025", 49, 48, 159.5, 2, 1357, 60, 22.61, 8.49, 15.9, 2, 0, "BOWL", 8, 0.97959183673469385, 2, 7.0221922305255111)
$row9 = @("E Malinga", "2024-2024", 1, 1, 1.5, 0, 26, 0, 100, 14.18, 100, 0, 0, "BOWL", 9, 1, 0, -14.678884700999999)

for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8[$i]
}

for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

# Copy formatting from the prior data row to the newly added rows
$ws.Range("A7:R7").Copy() | Out-Null
$ws.Range("A8:R9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the actual cell values (paste special copies formats only, but be safe)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8[$i]
}
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

# Update the selection to match the post-edit state
$ws.Range("I17").Select() | Out-Null
